$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: current (I [mA]) changes unit label to I [A]
$ws.Range("B1").Value2 = "I [A]"

# Convert B2:B20 current readings from mA to A (divide by 1000)
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 / 1000
}

# Update selection to H4 (matches final workbook selection state)
$ws.Range("H4").Select() | Out-Null
